# V1.0_beta prototype. Firwmare fixes
# Update the voltage regulator BOM line (row 11): the part changed from
# LM3940IMP-3.3/NOPB (SOT-223, JLCPCB C140319) to
# LT1086IM-3.3#TRPBF (JLCPCB C662441), and the footprint column is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the JLCPCB part # first, then the new component name, matching the
# order the strings were authored in (keeps shared-string table order
# consistent with the source edit).
$ws.Range("D11").Value = "C662441"
$ws.Range("A11").Value = "LT1086IM-3.3#TRPBF"
$ws.Range("C11").ClearContents()

# Leave the active selection on A11, as in the saved workbook.
[void]$ws.Range("A11").Select()
